$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")

# Row 19: confidence note replaced with "High" (note text moved down to row 21)
$ws.Range("F19").Value = "High"

# Row 21: second DFS binary tree problem logged
$ws.Range("B21").Value = "Binary Tree - DFS - LC75"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 872
$ws.Range("F21").Value = "Need to pick a med for tomorrow as well. Stil half a day behind. "

# Row 22: medium quantity bumped up
$ws.Range("D22").Value = 2

# Update the active selection
$ws.Range("G18").Select() | Out-Null
